$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.170.15"
$ws.Range("E2").Value = "  +1.71%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.189.09"
$ws.Range("E3").Value = "  +4.13%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.33"
$ws.Range("E5").Value = "  +3.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.53"
$ws.Range("E6").Value = "  +5.53%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.187.03"
$ws.Range("E8").Value = "  +4.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("E9").Value = "  +3.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("E10").Value = "  +5.95%  "

$ws.Range("E11").Value = "  +1.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.509"
$ws.Range("E12").Value = "  +2.40%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000275"
$ws.Range("E13").Value = "  +18.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.25"
$ws.Range("E14").Value = "  +7.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.711.64"
$ws.Range("E15").Value = "  +4.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.252.19"
$ws.Range("E16").Value = "  +1.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.191.51"
$ws.Range("E17").Value = "  +4.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.21"
$ws.Range("E18").Value = "  +6.17%  "

$ws.Range("E19").Value = "  +1.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "516.10"
$ws.Range("E20").Value = "  +7.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.95"
$ws.Range("E21").Value = "  +6.95%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.738"
$ws.Range("E22").Value = "  +7.87%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.44"
$ws.Range("E23").Value = "  +6.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.87"
$ws.Range("E24").Value = "  +3.81%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.54"
$ws.Range("E25").Value = "  +3.72%  "

$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.06"
$ws.Range("E27").Value = "  +10.86%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.93"
$ws.Range("E28").Value = "  +4.13%  "

$ws.Range("E29").Value = "  +7.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.24"
$ws.Range("E30").Value = "  +7.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.79"
$ws.Range("E31").Value = "  +13.98%  "

$ws.Range("E32").Value = "  +6.39%  "

$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("E34").Value = "  +8.79%  "

$ws.Range("E35").Value = "  +6.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.74"
$ws.Range("E36").Value = "  +1.49%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0903"
$ws.Range("E37").Value = "  +10.68%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "480.53"
$ws.Range("E38").Value = "  +6.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.17"
$ws.Range("E39").Value = "  +11.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0423"
$ws.Range("E40").Value = "  +2.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.69"
$ws.Range("E41").Value = "  +4.83%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.080.60"
$ws.Range("E42").Value = "  +1.99%  "

$ws.Range("E43").Value = "  +3.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.290"
$ws.Range("E44").Value = "  +8.74%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.44"
$ws.Range("E45").Value = "  +9.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.48"
$ws.Range("E46").Value = "  +5.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0606"
$ws.Range("E47").Value = "  +16.97%  "

$ws.Range("E48").Value = "  -0.08%  "

$ws.Range("E49").Value = "  +2.43%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.30"
$ws.Range("E50").Value = "  +10.17%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.34"
$ws.Range("E51").Value = "  +2.57%  "
